$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

# Title / headings (appears twice, identical replacement each time)
Replace-Text "Play Naughty Santa Slot Game for Free - Review" "Play Naughty Santa Slot for Free"
Replace-Text "Play Naughty Santa Slot Game for Free - Review" "Play Naughty Santa Slot for Free"

# "What we like" bullet list
Replace-Text "Fully compatible with mobile devices" "Mobile compatible for on-the-go gaming"
Replace-Text "Range of bonus features including free spins and multipliers" "Offers free play and real money options"
Replace-Text "High maximum win per spin of 57,600x your bet" "Exciting bonus features with free spins and multipliers"
Replace-Text "Can be played for free in demo mode" "High volatility with a maximum win of 57,600x"

# "What we don't like" bullet list
Replace-Text "High volatility may not be suitable for beginners" "May be too risky for inexperienced players"
Replace-Text "RTP range is quite wide (92.03% - 97.98%)" "Limited betting range with a maximum bet of €120.00"

# Closing meta description paragraph
Replace-Text "Read our review of Naughty Santa, a high volatility slot game with free play and real money options. Play now for free and discover bonus features." "Experience the excitement of Naughty Santa slot game with free play and bonus features."
